$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(208, 6).Value = 5069
$ws.Cells.Item(208, 7).Value = 14088.33333333333
$ws.Cells.Item(208, 8).Value = 321326.3333333333

$ws.Cells.Item(209, 6).Value = 5022.888888888889
$ws.Cells.Item(209, 7).Value = 13843
$ws.Cells.Item(209, 8).Value = 320897.6666666666

$ws.Cells.Item(210, 6).Value = 5038.296296296297
$ws.Cells.Item(210, 7).Value = 13900.22222222222
$ws.Cells.Item(210, 8).Value = 320238.4444444444

$ws.Cells.Item(211, 6).Value = 5043.395061728395
$ws.Cells.Item(211, 7).Value = 13943.85185185185
$ws.Cells.Item(211, 8).Value = 320820.8148148148

$ws.Cells.Item(212, 6).Value = 5034.860082304527
$ws.Cells.Item(212, 7).Value = 13895.69135802469
$ws.Cells.Item(212, 8).Value = 320652.3086419753

$ws.Cells.Item(213, 6).Value = 5038.85048010974
$ws.Cells.Item(213, 7).Value = 13913.25514403292
$ws.Cells.Item(213, 8).Value = 320570.5226337448

$ws.Cells.Item(214, 6).Value = 5039.035208047554
$ws.Cells.Item(214, 7).Value = 13917.59945130316
$ws.Cells.Item(214, 8).Value = 320681.2153635116

$ws.Cells.Item(215, 6).Value = 5037.581923487273
$ws.Cells.Item(215, 7).Value = 13908.84865112026
$ws.Cells.Item(215, 8).Value = 320634.6822130773

$ws.Cells.Item(216, 6).Value = 5038.489203881522
$ws.Cells.Item(216, 7).Value = 13913.23441548544
$ws.Cells.Item(216, 8).Value = 320628.8067367779

$ws.Cells.Item(217, 6).Value = 5038.368778472117
$ws.Cells.Item(217, 7).Value = 13913.22750596962
$ws.Cells.Item(217, 8).Value = 320648.2347711223

$ws.Cells.Item(218, 6).Value = 5038.146635280305
$ws.Cells.Item(218, 7).Value = 13911.77019085844
$ws.Cells.Item(218, 8).Value = 320637.2412403258

$ws.Cells.Item(219, 6).Value = 5038.334872544648
$ws.Cells.Item(219, 7).Value = 13912.74403743783
$ws.Cells.Item(219, 8).Value = 320638.0942494086

$ws.Cells.Item(220, 6).Value = 5038.28342876569
$ws.Cells.Item(220, 7).Value = 13912.58057808863
$ws.Cells.Item(220, 8).Value = 320641.1900869522

$ws.Cells.Item(221, 6).Value = 5038.254978863547
$ws.Cells.Item(221, 7).Value = 13912.36493546163
$ws.Cells.Item(221, 8).Value = 320638.8418588956

$ws.Cells.Item(222, 6).Value = 5038.291093391295
$ws.Cells.Item(222, 7).Value = 13912.5631836627
$ws.Cells.Item(222, 8).Value = 320639.3753984188

$ws.Cells.Item(223, 6).Value = 5038.276500340177
$ws.Cells.Item(223, 7).Value = 13912.50289907099
$ws.Cells.Item(223, 8).Value = 320639.8024480889

$ws.Cells.Item(224, 6).Value = 5038.274190865006
$ws.Cells.Item(224, 7).Value = 13912.47700606511
$ws.Cells.Item(224, 8).Value = 320639.3399018011

$ws.Cells.Item(225, 6).Value = 5038.280594865492
$ws.Cells.Item(225, 7).Value = 13912.51436293293
$ws.Cells.Item(225, 8).Value = 320639.5059161029

$ws.Cells.Item(226, 6).Value = 5038.277095356892
$ws.Cells.Item(226, 7).Value = 13912.49808935634
$ws.Cells.Item(226, 8).Value = 320639.5494219976

$ws.Cells.Item(227, 6).Value = 5038.277293695796
$ws.Cells.Item(227, 7).Value = 13912.49648611813
$ws.Cells.Item(227, 8).Value = 320639.4650799672

$ws.Cells.Item(228, 6).Value = 5038.278327972726
$ws.Cells.Item(228, 7).Value = 13912.50297946913
$ws.Cells.Item(228, 8).Value = 320639.5068060226

$ws.Cells.Item(229, 6).Value = 5038.277572341804
$ws.Cells.Item(229, 7).Value = 13912.4991849812
$ws.Cells.Item(229, 8).Value = 320639.5071026625

$ws.Cells.Item(230, 6).Value = 5038.277731336776
$ws.Cells.Item(230, 7).Value = 13912.49955018949
$ws.Cells.Item(230, 8).Value = 320639.4929962175

$ws.Cells.Item(231, 6).Value = 5038.277877217102
$ws.Cells.Item(231, 7).Value = 13912.50057154661
$ws.Cells.Item(231, 8).Value = 320639.5023016342

$ws.Cells.Item(232, 6).Value = 5038.277726965228
$ws.Cells.Item(232, 7).Value = 13912.49976890576
$ws.Cells.Item(232, 8).Value = 320639.5008001713

$ws.Cells.Item(233, 6).Value = 5038.277778506369
$ws.Cells.Item(233, 7).Value = 13912.49996354729
$ws.Cells.Item(233, 8).Value = 320639.498699341

$ws.Cells.Item(234, 6).Value = 5038.277794229566
$ws.Cells.Item(234, 7).Value = 13912.50010133322
$ws.Cells.Item(234, 8).Value = 320639.5006003822

$ws.Cells.Item(235, 6).Value = 5038.277766567055
$ws.Cells.Item(235, 7).Value = 13912.49994459542
$ws.Cells.Item(235, 8).Value = 320639.5000332982

$ws.Cells.Item(236, 6).Value = 5038.277779767664
$ws.Cells.Item(236, 7).Value = 13912.50000315864
$ws.Cells.Item(236, 8).Value = 320639.4997776737

$ws.Cells.Item(237, 6).Value = 5038.277780188095
$ws.Cells.Item(237, 7).Value = 13912.50001636243
$ws.Cells.Item(237, 8).Value = 320639.500137118

$ws.Cells.Item(238, 6).Value = 5038.277775507604
$ws.Cells.Item(238, 7).Value = 13912.49998803883
$ws.Cells.Item(238, 8).Value = 320639.4999826967

$ws.Cells.Item(239, 6).Value = 5038.277778487787
$ws.Cells.Item(239, 7).Value = 13912.50000251997
$ws.Cells.Item(239, 8).Value = 320639.4999658295

$ws.Cells.Item(240, 6).Value = 5038.277778061162
$ws.Cells.Item(240, 7).Value = 13912.50000230708
$ws.Cells.Item(240, 8).Value = 320639.5000285481

$ws.Cells.Item(241, 6).Value = 5038.277777352185
$ws.Cells.Item(241, 7).Value = 13912.49999762196
$ws.Cells.Item(241, 8).Value = 320639.4999923581

$ws.Cells.Item(242, 6).Value = 5038.277777967045
$ws.Cells.Item(242, 7).Value = 13912.50000081633
$ws.Cells.Item(242, 8).Value = 320639.4999955786

$ws.Cells.Item(243, 6).Value = 5038.277777793464
$ws.Cells.Item(243, 7).Value = 13912.50000024846
$ws.Cells.Item(243, 8).Value = 320639.5000054949

$ws.Cells.Item(244, 6).Value = 5038.277777704231
$ws.Cells.Item(244, 7).Value = 13912.49999956225
$ws.Cells.Item(244, 8).Value = 320639.4999978105

$ws.Cells.Item(245, 6).Value = 5038.277777821581
$ws.Cells.Item(245, 7).Value = 13912.50000020901
$ws.Cells.Item(245, 8).Value = 320639.499999628

$ws.Cells.Item(246, 6).Value = 5038.277777773092
$ws.Cells.Item(246, 7).Value = 13912.50000000657
$ws.Cells.Item(246, 8).Value = 320639.5000009778

$ws.Cells.Item(247, 6).Value = 5038.277777766301
$ws.Cells.Item(247, 7).Value = 13912.49999992595
$ws.Cells.Item(247, 8).Value = 320639.4999994721

$ws.Cells.Item(248, 6).Value = 5038.277777786992
$ws.Cells.Item(248, 7).Value = 13912.50000004718
$ws.Cells.Item(248, 8).Value = 320639.500000026

$ws.Cells.Item(249, 6).Value = 5038.277777775462
$ws.Cells.Item(249, 7).Value = 13912.49999999323
$ws.Cells.Item(249, 8).Value = 320639.5000001586
